# Apply weekly fruit/vegetable price updates to the Papaya sheet.
# The edit re-distributes the per-row data block (Fecha, Calidad, Volumen,
# Precio minimo/maximo/promedio, Unidad de comercializacion, Precio $/Kg,
# Kg/unidad) across rows 2,3,7,8,9,10,11,12,13,14,15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowData {
    param($Row, $Fecha, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $Unidad, $PrecioKg, $KgUnidad)

    $ws.Cells.Item($Row, 4).Value2 = $Fecha         # D - Fecha
    $ws.Cells.Item($Row, 12).Value = $Calidad       # L - Calidad
    $ws.Cells.Item($Row, 13).Value2 = $Volumen      # M - Volumen
    $ws.Cells.Item($Row, 14).Value2 = $PrecioMin    # N - Precio minimo
    $ws.Cells.Item($Row, 15).Value2 = $PrecioMax    # O - Precio maximo
    $ws.Cells.Item($Row, 16).Value2 = $PrecioProm   # P - Precio promedio ponderado
    $ws.Cells.Item($Row, 17).Value = $Unidad        # Q - Unidad de comercializacion
    $ws.Cells.Item($Row, 19).Value2 = $PrecioKg     # S - Precio $/Kg
    $ws.Cells.Item($Row, 20).Value2 = $KgUnidad     # T - Kg / unidad
}

Set-RowData 2  44371 "Primera" 20  1800  1800  1800  "$/kilo (en caja de 15 kilos)" 1800 1
Set-RowData 3  44371 "Segunda" 30  1200  1200  1200  "$/kilo (en caja de 15 kilos)" 1200 1
Set-RowData 7  44880 "Primera" 200 20000 20000 20000 "$/bandeja 10 kilos"           2000 10
Set-RowData 8  44880 "Segunda" 180 15000 15000 15000 "$/bandeja 10 kilos"           1500 10
Set-RowData 9  44343 "Primera" 20  1700  1700  1700  "$/kilo (en caja de 15 kilos)" 1700 1
Set-RowData 10 44400 "Primera" 25  1500  1500  1500  "$/kilo (en caja de 15 kilos)" 1500 1
Set-RowData 11 44904 "Primera" 45  15000 15000 15000 "$/bandeja 10 kilos"           1500 10
Set-RowData 12 44904 "Segunda" 60  10000 10000 10000 "$/bandeja 10 kilos"           1000 10
Set-RowData 13 44336 "Primera" 10  1500  1500  1500  "$/kilo (en caja de 15 kilos)" 1500 1
Set-RowData 14 44292 "Primera" 50  14000 14000 14000 "$/bandeja 10 kilos"           1400 10
Set-RowData 15 44309 "Primera" 10  1600  1600  1600  "$/kilo (en caja de 15 kilos)" 1600 1
